# Edit script: add "2022-Q4" sheet with fund holding data, and update
# the "总计" (summary) sheet with the new quarter's totals.

$wb = $excel.ActiveWorkbook

# Helper: force a numeric-looking string to be stored as text (preserve
# leading/trailing zeros, e.g. fund codes "008969" or ratios "3.10").
function Set-TextCell($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

# ---------------------------------------------------------------------
# 1. Insert a brand new worksheet named "2022-Q4" right after "总计",
#    pushing every other quarter sheet down by one position.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$ns = $wb.Worksheets.Add($null, $summary)
$ns.Name = "2022-Q4"

# Copy the header / index-column formatting (style "s=2", bold + border)
# from the summary sheet so the new sheet matches the look of its
# siblings.
$summary.Range("B1:D1").Copy()
$ns.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$ns.Range("A2:A16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Populate the "2022-Q4" sheet with the fund holding data.
# ---------------------------------------------------------------------
$ns.Range("B1").Value = "基金代码"
$ns.Range("C1").Value = "基金名称"
$ns.Range("D1").Value = "基金规模"
$ns.Range("E1").Value = "股票总仓位"
$ns.Range("F1").Value = "仓位占比"
$ns.Range("G1").Value = "持有市值(亿元)"
$ns.Range("H1").Value = "仓位排名"
$ns.Range("A2").Value = 0
Set-TextCell $ns "B2" "008969"
$ns.Range("C2").Value = "睿远均衡价值三年持有期混合A"
Set-TextCell $ns "D2" "134.09"
Set-TextCell $ns "E2" "91.04"
Set-TextCell $ns "F2" "2.94"
Set-TextCell $ns "G2" "3.9422"
$ns.Range("H2").Value = 10
$ns.Range("A3").Value = 1
Set-TextCell $ns "B3" "010902"
$ns.Range("C3").Value = "博时成长领航灵活配置混合A"
Set-TextCell $ns "D3" "48.67"
Set-TextCell $ns "E3" "89.03"
Set-TextCell $ns "F3" "5.77"
Set-TextCell $ns "G3" "2.8083"
$ns.Range("H3").Value = 4
$ns.Range("A4").Value = 2
Set-TextCell $ns "B4" "513060"
$ns.Range("C4").Value = "博时恒生医疗保健ETF（QDII）"
Set-TextCell $ns "D4" "69.51"
Set-TextCell $ns "E4" "99.65"
Set-TextCell $ns "F4" "3.19"
Set-TextCell $ns "G4" "2.2174"
$ns.Range("H4").Value = 8
$ns.Range("A5").Value = 3
Set-TextCell $ns "B5" "008970"
$ns.Range("C5").Value = "睿远均衡价值三年持有期混合C"
Set-TextCell $ns "D5" "16.67"
Set-TextCell $ns "E5" "91.04"
Set-TextCell $ns "F5" "2.94"
Set-TextCell $ns "G5" "0.4901"
$ns.Range("H5").Value = 10
$ns.Range("A6").Value = 4
Set-TextCell $ns "B6" "010903"
$ns.Range("C6").Value = "博时成长领航灵活配置混合C"
Set-TextCell $ns "D6" "7.15"
Set-TextCell $ns "E6" "89.03"
Set-TextCell $ns "F6" "5.77"
Set-TextCell $ns "G6" "0.4126"
$ns.Range("H6").Value = 4
$ns.Range("A7").Value = 5
Set-TextCell $ns "B7" "159892"
$ns.Range("C7").Value = "华夏恒生香港上市生物科技ETF（QDII）"
Set-TextCell $ns "D7" "5.36"
Set-TextCell $ns "E7" "99.47"
Set-TextCell $ns "F7" "3.52"
Set-TextCell $ns "G7" "0.1887"
$ns.Range("H7").Value = 8
$ns.Range("A8").Value = 6
Set-TextCell $ns "B8" "513700"
$ns.Range("C8").Value = "鹏华中证港股通医药卫生综合ETF"
Set-TextCell $ns "D8" "4.34"
Set-TextCell $ns "E8" "95.37"
Set-TextCell $ns "F8" "2.69"
Set-TextCell $ns "G8" "0.1167"
$ns.Range("H8").Value = 8
$ns.Range("A9").Value = 7
Set-TextCell $ns "B9" "513200"
$ns.Range("C9").Value = "易方达中证港股通医药卫生综合ETF"
Set-TextCell $ns "D9" "1.69"
Set-TextCell $ns "E9" "95.67"
Set-TextCell $ns "F9" "2.76"
Set-TextCell $ns "G9" "0.0466"
$ns.Range("H9").Value = 8
$ns.Range("A10").Value = 8
Set-TextCell $ns "B10" "513280"
$ns.Range("C10").Value = "汇添富恒生香港上市生物科技ETF（QDII）"
Set-TextCell $ns "D10" "1.35"
Set-TextCell $ns "E10" "94.55"
Set-TextCell $ns "F10" "3.38"
Set-TextCell $ns "G10" "0.0456"
$ns.Range("H10").Value = 8
$ns.Range("A11").Value = 9
Set-TextCell $ns "B11" "159776"
$ns.Range("C11").Value = "银华中证港股通医药卫生综合ETF"
Set-TextCell $ns "D11" "0.81"
Set-TextCell $ns "E11" "93.98"
Set-TextCell $ns "F11" "2.65"
Set-TextCell $ns "G11" "0.0215"
$ns.Range("H11").Value = 8
$ns.Range("A12").Value = 10
Set-TextCell $ns "B12" "159718"
$ns.Range("C12").Value = "平安中证港股通医药卫生综合ETF"
Set-TextCell $ns "D12" "0.67"
Set-TextCell $ns "E12" "94.63"
Set-TextCell $ns "F12" "2.65"
Set-TextCell $ns "G12" "0.0178"
$ns.Range("H12").Value = 8
$ns.Range("A13").Value = 11
Set-TextCell $ns "B13" "008861"
$ns.Range("C13").Value = "西部利得港股通新机遇灵活配置混合A"
Set-TextCell $ns "D13" "0.25"
Set-TextCell $ns "E13" "87.69"
Set-TextCell $ns "F13" "3.39"
Set-TextCell $ns "G13" "0.0085"
$ns.Range("H13").Value = 9
$ns.Range("A14").Value = 12
Set-TextCell $ns "B14" "010093"
$ns.Range("C14").Value = "西部利得港股通新机遇灵活配置混合C"
Set-TextCell $ns "D14" "0.12"
Set-TextCell $ns "E14" "87.69"
Set-TextCell $ns "F14" "3.39"
Set-TextCell $ns "G14" "0.0041"
$ns.Range("H14").Value = 9
$ns.Range("A15").Value = 13
Set-TextCell $ns "B15" "012315"
$ns.Range("C15").Value = "创金合信港股通成长股票A"
Set-TextCell $ns "D15" "0.12"
Set-TextCell $ns "E15" "89.18"
Set-TextCell $ns "F15" "3.10"
Set-TextCell $ns "G15" "0.0037"
$ns.Range("H15").Value = 9
$ns.Range("A16").Value = 14
Set-TextCell $ns "B16" "012316"
$ns.Range("C16").Value = "创金合信港股通成长股票C"
Set-TextCell $ns "D16" "0.11"
Set-TextCell $ns "E16" "89.18"
Set-TextCell $ns "F16" "3.10"
Set-TextCell $ns "G16" "0.0034"
$ns.Range("H16").Value = 9

# ---------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: the new quarter becomes the
#    first data row, and every previous quarter row shifts down by one.
# ---------------------------------------------------------------------
# Row 9 is brand new (previously the sheet only had rows 1-8); copy the
# index-column formatting (style "s=2") from row 8 onto it.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 10.33

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 11
$summary.Range("D3").Value = 4.61

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 1.52

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.6

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 9
$summary.Range("D6").Value = 7.13

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 21
$summary.Range("D7").Value = 32.14

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 15
$summary.Range("D8").Value = 28.33

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 11
$summary.Range("D9").Value = 14.43
